{"js": "// Grammar fix: \"Ils devaient des demander\" -> \"Ils devaient se demander\"\n// (and the Word \"_GoBack\" last-edit bookmark moves along with the edit,\n// from its old spot near \" Alors, q\" to right after the newly typed \"se\").\nconst body = context.document.body;\n\n// 1) Find the erroneous word \"des\" in its sentence context and replace it\n//    with \"se\", leaving the rest of the sentence untouched.\nconst hits = body.search(\"devaient des demander\", { matchCase: true, matchWholeWord: false });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length === 0) {\n  throw new Error(\"Could not find the target phrase 'devaient des demander'\");\n}\n\nconst phrase = hits.items[0];\n// Split the hit into its three whitespace-delimited words so we can target\n// \"des\" precisely without disturbing the surrounding \" demander\" text.\nconst words = phrase.split([\" \"], false, true, true);\nwords.load(\"items/text\");\nawait context.sync();\n\nconst desWord = words.items[1]; // \"des\"\ndesWord.insertText(\"se\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) Move the \"_GoBack\" bookmark: remove it from its old location and drop\n//    it right after the freshly corrected \"se\", matching where Word leaves\n//    its \"last edit\" marker after typing at that spot.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst fixedHits = body.search(\"Ils devaient se demander\", { matchCase: true });\nfixedHits.load(\"items\");\nawait context.sync();\n\nconst fixedSentence = fixedHits.items[0];\nconst fixedWords = fixedSentence.split([\" \"], false, true, true);\nfixedWords.load(\"items/text\");\nawait context.sync();\n\nconst seWord = fixedWords.items[2]; // \"se\"\nconst afterSe = seWord.getRange(\"End\");\nafterSe.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Grammar fix: \"Ils devaient des demander\" -> \"Ils devaient se demander\"\n# (and the Word \"_GoBack\" last-edit bookmark moves along with the edit,\n# from its old spot near \" Alors, q\" to right after the newly typed \"se\").\n$d = $word.ActiveDocument\n\n# 1) Replace the erroneous \"des\" with \"se\" in its sentence context, leaving\n#    the rest of the sentence (\"... demander ce qui poussait ...\") intact.\n$findRange = $d.Content\n$replaced = $findRange.Find.Execute(\"devaient des demander\", $false, $false, $false, $false, $false, $true, 1, $false, \"devaient se demander\", 2)\nif (-not $replaced) {\n    throw \"Could not find the target phrase 'devaient des demander'\"\n}\n\n# 2) Move the \"_GoBack\" bookmark: drop it from its old location and re-add\n#    it right after the freshly corrected \"se\", matching where Word leaves\n#    its \"last edit\" marker after typing at that spot.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$seRange = $d.Content\n$seRange.Find.Execute(\"Ils devaient se\")\n$bookmarkRange = $d.Range($seRange.End, $seRange.End)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n"}
